$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1019.8
$ws.Range("I18").Value = 774.75
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 774.75
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -490.75
$ws.Range("N18").Value = -2568
$ws.Range("H28").Value = 117066.11
$ws.Range("I28").Value = 146869.86
$ws.Range("J28").Value = 12753
$ws.Range("K28").Value = 146869.86
$ws.Range("L28").Value = 12753
$ws.Range("M28").Value = -146384.86
$ws.Range("N28").Value = -13723
$ws.Range("H33").Value = 364.3158
$ws.Range("I33").Value = 170.875
$ws.Range("J33").Value = 1396
$ws.Range("K33").Value = 170.875
$ws.Range("L33").Value = 1396
$ws.Range("M33").Value = 58.125
$ws.Range("N33").Value = -1854
$ws.Range("H86").Value = 7106.875
$ws.Range("I86").Value = 7120
$ws.Range("J86").Value = 7085
$ws.Range("K86").Value = 7120
$ws.Range("L86").Value = 7085
$ws.Range("M86").Value = -5997
$ws.Range("N86").Value = -9331
$ws.Range("H89").Value = 7106.875
$ws.Range("I89").Value = 7120
$ws.Range("J89").Value = 7085
$ws.Range("K89").Value = 35600
$ws.Range("L89").Value = 35425
$ws.Range("M89").Value = -29984
$ws.Range("N89").Value = -46657
$ws.Range("H113").Value = 14181.6
$ws.Range("I113").Value = 8000
$ws.Range("J113").Value = 15727
$ws.Range("K113").Value = 8000
$ws.Range("L113").Value = 15727
$ws.Range("M113").Value = -4746
$ws.Range("N113").Value = -22235
$ws.Range("H133").Value = 61982.668
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 61982.668
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 61982.668
$ws.Range("N133").Value = -72102.66800000001
$ws.Range("H137").Value = 3060.5833
$ws.Range("I137").Value = 1679.75
$ws.Range("J137").Value = 3336.75
$ws.Range("K137").Value = 5039.25
$ws.Range("L137").Value = 10010.25
$ws.Range("M137").Value = -2489.25
$ws.Range("N137").Value = -15110.25
$ws.Range("H138").Value = 3761.6052
$ws.Range("I138").Value = 11000
$ws.Range("J138").Value = 3565.973
$ws.Range("K138").Value = 33000
$ws.Range("L138").Value = 10697.919
$ws.Range("M138").Value = -27860
$ws.Range("N138").Value = -20977.919
$ws.Range("H141").Value = 7214.45
$ws.Range("I141").Value = 4498.8
$ws.Range("J141").Value = 8119.6665
$ws.Range("K141").Value = 13496.4
$ws.Range("L141").Value = 24358.9995
$ws.Range("M141").Value = -8316.400000000001
$ws.Range("N141").Value = -34718.99950000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 111116600
$ws.Range("I45").Value = 250001500
$ws.Range("J45").Value = 8672.200000000001
$ws.Range("K45").Value = 250001500
$ws.Range("L45").Value = 8672.200000000001
$ws.Range("M45").Value = -250001123
$ws.Range("N45").Value = -9426.200000000001
$ws.Range("H61").Value = 3554.9412
$ws.Range("I61").Value = 3554.9412
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3554.9412
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3342.9412
$ws.Range("H74").Value = 55558880
$ws.Range("I74").Value = 83337320
$ws.Range("J74").Value = 2007
$ws.Range("K74").Value = 83337320
$ws.Range("L74").Value = 2007
$ws.Range("M74").Value = -83336446
$ws.Range("N74").Value = -3755
$ws.Range("H77").Value = 55558880
$ws.Range("I77").Value = 83337320
$ws.Range("J77").Value = 2007
$ws.Range("K77").Value = 416686600
$ws.Range("L77").Value = 10035
$ws.Range("M77").Value = -416682232
$ws.Range("N77").Value = -18771
$ws.Range("H110").Value = 9516.4
$ws.Range("I110").Value = 8123.2856
$ws.Range("J110").Value = 12767
$ws.Range("K110").Value = 8123.2856
$ws.Range("L110").Value = 12767
$ws.Range("M110").Value = -6078.2856
$ws.Range("N110").Value = -16857
$ws.Range("H119").Value = 37336
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 37336
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 37336
$ws.Range("N119").Value = -47012
$ws.Range("H122").Value = 6750
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 7666.6665
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 22999.9995
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -27899.9995
$ws.Range("H132").Value = 7585.364
$ws.Range("I132").Value = 2272.3333
$ws.Range("J132").Value = 9577.75
$ws.Range("K132").Value = 6816.999899999999
$ws.Range("L132").Value = 28733.25
$ws.Range("M132").Value = -4286.999899999999
$ws.Range("N132").Value = -33793.25
$ws.Range("H136").Value = 3554.9412
$ws.Range("I136").Value = 3554.9412
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10664.8236
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8114.8236

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 29244.416
$ws.Range("I50").Value = 27491.5
$ws.Range("J50").Value = 29595
$ws.Range("K50").Value = 27491.5
$ws.Range("L50").Value = 29595
$ws.Range("M50").Value = -26866.5
$ws.Range("N50").Value = -30845
$ws.Range("H51").Value = 30000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 30000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 30000
$ws.Range("N51").Value = -31472
$ws.Range("M51").ClearContents() | Out-Null
$ws.Range("H60").Value = 14150
$ws.Range("I60").Value = 6225
$ws.Range("J60").Value = 30000
$ws.Range("K60").Value = 6225
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = -5714
$ws.Range("N60").Value = -31022
$ws.Range("H61").Value = 30000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 30000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 30000
$ws.Range("N61").Value = -30696
$ws.Range("M61").ClearContents() | Out-Null
$ws.Range("H86").Value = 11099.7
$ws.Range("I86").Value = 5873.5
$ws.Range("J86").Value = 14583.833
$ws.Range("K86").Value = 5873.5
$ws.Range("L86").Value = 14583.833
$ws.Range("M86").Value = -4750.5
$ws.Range("N86").Value = -16829.833
$ws.Range("H89").Value = 11099.7
$ws.Range("I89").Value = 5873.5
$ws.Range("J89").Value = 14583.833
$ws.Range("K89").Value = 29367.5
$ws.Range("L89").Value = 72919.16500000001
$ws.Range("M89").Value = -23751.5
$ws.Range("N89").Value = -84151.16500000001
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents() | Out-Null
$ws.Range("N99").ClearContents() | Out-Null
$ws.Range("H105").Value = 1980.5
$ws.Range("I105").Value = 1002.5333
$ws.Range("J105").Value = 6870.3335
$ws.Range("K105").Value = 1002.5333
$ws.Range("L105").Value = 6870.3335
$ws.Range("M105").Value = 744.4666999999999
$ws.Range("N105").Value = -10364.3335
$ws.Range("H107").Value = 853.2069
$ws.Range("I107").Value = 746.6
$ws.Range("J107").Value = 1519.5
$ws.Range("K107").Value = 746.6
$ws.Range("L107").Value = 1519.5
$ws.Range("M107").Value = 1173.4
$ws.Range("N107").Value = -5359.5
$ws.Range("H122").Value = 10353.8
$ws.Range("I122").Value = 5299.6665
$ws.Range("J122").Value = 17935
$ws.Range("K122").Value = 15898.9995
$ws.Range("L122").Value = 53805
$ws.Range("M122").Value = -13448.9995
$ws.Range("N122").Value = -58705
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents() | Out-Null
$ws.Range("N126").ClearContents() | Out-Null
$ws.Range("H132").Value = 7477.56
$ws.Range("I132").Value = 7002.5264
$ws.Range("J132").Value = 8981.833000000001
$ws.Range("K132").Value = 21007.5792
$ws.Range("L132").Value = 26945.499
$ws.Range("M132").Value = -18477.5792
$ws.Range("N132").Value = -32005.499
$ws.Range("H134").Value = 2350.4285
$ws.Range("I134").Value = 1730.6923
$ws.Range("J134").Value = 10407
$ws.Range("K134").Value = 5192.0769
$ws.Range("L134").Value = 31221
$ws.Range("M134").Value = -2657.0769
$ws.Range("N134").Value = -36291

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4146.9375
$ws.Range("I132").Value = 2809.9
$ws.Range("J132").Value = 6375.3335
$ws.Range("K132").Value = 25289.1
$ws.Range("L132").Value = 57378.0015
$ws.Range("M132").Value = -22759.1
$ws.Range("N132").Value = -62438.0015

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14207.772
$ws.Range("I70").Value = 5184.0713
$ws.Range("J70").Value = 29999.25
$ws.Range("K70").Value = 5184.0713
$ws.Range("L70").Value = 29999.25
$ws.Range("M70").Value = -4914.0713
$ws.Range("N70").Value = -30539.25
$ws.Range("H73").Value = 14207.772
$ws.Range("I73").Value = 5184.0713
$ws.Range("J73").Value = 29999.25
$ws.Range("K73").Value = 5184.0713
$ws.Range("L73").Value = 29999.25
$ws.Range("M73").Value = -4248.0713
$ws.Range("N73").Value = -31871.25
$ws.Range("H102").Value = 3626.7778
$ws.Range("I102").Value = 1558.25
$ws.Range("J102").Value = 5281.6
$ws.Range("K102").Value = 1558.25
$ws.Range("L102").Value = 5281.6
$ws.Range("M102").Value = 63.75
$ws.Range("N102").Value = -8525.6
$ws.Range("H126").Value = 4666.6665
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 155572.86
$ws.Range("I132").Value = 342332.34
$ws.Range("J132").Value = 15503.25
$ws.Range("K132").Value = 1026997.02
$ws.Range("L132").Value = 46509.75
$ws.Range("M132").Value = -1024467.02
$ws.Range("N132").Value = -51569.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6575.8887
$ws.Range("I40").Value = 8321.5
$ws.Range("J40").Value = 4393.875
$ws.Range("K40").Value = 8321.5
$ws.Range("L40").Value = 4393.875
$ws.Range("M40").Value = -8185.5
$ws.Range("N40").Value = -4665.875
$ws.Range("H132").Value = 4257.8945
$ws.Range("I132").Value = 1641.75
$ws.Range("J132").Value = 8742.714
$ws.Range("K132").Value = 4925.25
$ws.Range("L132").Value = 26228.142
$ws.Range("M132").Value = -2395.25
$ws.Range("N132").Value = -31288.142

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3526.5293
$ws.Range("I132").Value = 3370.0667
$ws.Range("J132").Value = 4700
$ws.Range("K132").Value = 10110.2001
$ws.Range("L132").Value = 14100
$ws.Range("M132").Value = -7580.2001
$ws.Range("N132").Value = -19160
$ws.Range("H136").Value = 8727.727999999999
$ws.Range("I136").Value = 6444.4443
$ws.Range("J136").Value = 19002.5
$ws.Range("K136").Value = 19333.3329
$ws.Range("L136").Value = 57007.5
$ws.Range("M136").Value = -16783.3329
$ws.Range("N136").Value = -62107.5
